# "Generate Report for Archive" — refresh the localization-status report:
#   - the handoff status for the sample row moved from "Ready for handoff"
#     to "In Translation" on every sheet that surfaces it
#   - the Status column (narrower now that "In Translation" is shorter than
#     "Ready for handoff") is re-sized to fit the new text

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Update the status text everywhere it appears.
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value     = "In Translation"
$dede.Range("C2").Value     = "In Translation"

# Re-fit the now-narrower status columns to the shorter text.
$overview.Range("E1").ColumnWidth = 12.5
$overview.Range("F1").ColumnWidth = 12.5
$zhcn.Range("C1").ColumnWidth     = 12.5
$dede.Range("C1").ColumnWidth     = 12.5
